# Natmi following Dr Hou advice
# Adds an "ECs" sending/target cluster to the Lgi2-Adam11 LR-pair table,
# replacing the old 2x2 (FAPs/sCs) grid with a full 3x3 (ECs/FAPs/sCs) grid
# and refreshed edge-expression statistics for every cluster pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: ECs -> ECs (Lgi2/Adam11)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Lgi2"
$ws.Range("C2").Value = "Adam11"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1.0
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.058936
$ws.Range("H2").Value = 0.176808
$ws.Range("I2").Value = 0.003640699631737656
$ws.Range("J2").Value = 0.003640699631737656
$ws.Range("K2").Value = 1.0
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.07360033333333334
$ws.Range("N2").Value = 0.220801
$ws.Range("O2").Value = 0.02873450582079328
$ws.Range("P2").Value = 0.02873450582079327
$ws.Range("Q2").Value = 0.004337709245333333
$ws.Range("R2").Value = 0.039039383208
$ws.Range("S2").Value = 0.0001046137047599256
$ws.Range("T2").Value = 0.0001046137047599256

# Row 3: ECs -> FAPs (Lgi2/Adam11)
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Lgi2"
$ws.Range("C3").Value = "Adam11"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1.0
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.058936
$ws.Range("H3").Value = 0.176808
$ws.Range("I3").Value = 0.003640699631737656
$ws.Range("J3").Value = 0.003640699631737656
$ws.Range("K3").Value = 3.0
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 0.4213883333333333
$ws.Range("N3").Value = 1.264165
$ws.Range("O3").Value = 0.1645153624799848
$ws.Range("P3").Value = 0.1645153624799848
$ws.Range("Q3").Value = 0.02483494281333333
$ws.Range("R3").Value = 0.22351448532
$ws.Range("S3").Value = 0.0005989510195960678
$ws.Range("T3").Value = 0.0005989510195960678

# Row 4: ECs -> sCs (Lgi2/Adam11)
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Lgi2"
$ws.Range("C4").Value = "Adam11"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1.0
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.058936
$ws.Range("H4").Value = 0.176808
$ws.Range("I4").Value = 0.003640699631737656
$ws.Range("J4").Value = 0.003640699631737656
$ws.Range("K4").Value = 3.0
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 2.066403333333334
$ws.Range("N4").Value = 6.199210000000001
$ws.Range("O4").Value = 0.8067501316992219
$ws.Range("P4").Value = 0.8067501316992219
$ws.Range("Q4").Value = 0.1217855468533333
$ws.Range("R4").Value = 1.09606992168
$ws.Range("S4").Value = 0.002937134907381663
$ws.Range("T4").Value = 0.002937134907381663

# Row 5: FAPs -> ECs (Lgi2/Adam11)
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Lgi2"
$ws.Range("C5").Value = "Adam11"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 14.32797533333333
$ws.Range("H5").Value = 42.983926
$ws.Range("I5").Value = 0.8850932285803735
$ws.Range("J5").Value = 0.8850932285803734
$ws.Range("K5").Value = 1.0
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.07360033333333334
$ws.Range("N5").Value = 0.220801
$ws.Range("O5").Value = 0.02873450582079328
$ws.Range("P5").Value = 0.02873450582079327
$ws.Range("Q5").Value = 1.054543760525111
$ws.Range("R5").Value = 9.490893844726
$ws.Range("S5").Value = 0.02543271652858746
$ws.Range("T5").Value = 0.02543271652858745

# Row 6: FAPs -> FAPs (Lgi2/Adam11)
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Lgi2"
$ws.Range("C6").Value = "Adam11"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 14.32797533333333
$ws.Range("H6").Value = 42.983926
$ws.Range("I6").Value = 0.8850932285803735
$ws.Range("J6").Value = 0.8850932285803734
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 0.4213883333333333
$ws.Range("N6").Value = 1.264165
$ws.Range("O6").Value = 0.1645153624799848
$ws.Range("P6").Value = 0.1645153624799848
$ws.Range("Q6").Value = 6.037641645754444
$ws.Range("R6").Value = 54.33877481179
$ws.Range("S6").Value = 0.1456114333284802
$ws.Range("T6").Value = 0.1456114333284802

# Row 7: FAPs -> sCs (Lgi2/Adam11)
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Lgi2"
$ws.Range("C7").Value = "Adam11"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 14.32797533333333
$ws.Range("H7").Value = 42.983926
$ws.Range("I7").Value = 0.8850932285803735
$ws.Range("J7").Value = 0.8850932285803734
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 2.066403333333334
$ws.Range("N7").Value = 6.199210000000001
$ws.Range("O7").Value = 0.8067501316992219
$ws.Range("P7").Value = 0.8067501316992219
$ws.Range("Q7").Value = 29.60737598871778
$ws.Range("R7").Value = 266.46638389846
$ws.Range("S7").Value = 0.7140490787233059
$ws.Range("T7").Value = 0.7140490787233058

# Row 8: sCs -> ECs (Lgi2/Adam11)
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Lgi2"
$ws.Range("C8").Value = "Adam11"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3.0
$ws.Range("F8").Value = 1.0
$ws.Range("G8").Value = 1.801186
$ws.Range("H8").Value = 5.403558
$ws.Range("I8").Value = 0.111266071787889
$ws.Range("J8").Value = 0.1112660717878889
$ws.Range("K8").Value = 1.0
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.07360033333333334
$ws.Range("N8").Value = 0.220801
$ws.Range("O8").Value = 0.02873450582079328
$ws.Range("P8").Value = 0.02873450582079327
$ws.Range("Q8").Value = 0.1325678899953334
$ws.Range("R8").Value = 1.193111009958
$ws.Range("S8").Value = 0.003197175587445898
$ws.Range("T8").Value = 0.003197175587445897

# Row 9: sCs -> FAPs (Lgi2/Adam11)
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Lgi2"
$ws.Range("C9").Value = "Adam11"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3.0
$ws.Range("F9").Value = 1.0
$ws.Range("G9").Value = 1.801186
$ws.Range("H9").Value = 5.403558
$ws.Range("I9").Value = 0.111266071787889
$ws.Range("J9").Value = 0.1112660717878889
$ws.Range("K9").Value = 3.0
$ws.Range("L9").Value = 1.0
$ws.Range("M9").Value = 0.4213883333333333
$ws.Range("N9").Value = 1.264165
$ws.Range("O9").Value = 0.1645153624799848
$ws.Range("P9").Value = 0.1645153624799848
$ws.Range("Q9").Value = 0.7589987665633333
$ws.Range("R9").Value = 6.83098889907
$ws.Range("S9").Value = 0.01830497813190857
$ws.Range("T9").Value = 0.01830497813190856

# Row 10: sCs -> sCs (Lgi2/Adam11)
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Lgi2"
$ws.Range("C10").Value = "Adam11"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3.0
$ws.Range("F10").Value = 1.0
$ws.Range("G10").Value = 1.801186
$ws.Range("H10").Value = 5.403558
$ws.Range("I10").Value = 0.111266071787889
$ws.Range("J10").Value = 0.1112660717878889
$ws.Range("K10").Value = 3.0
$ws.Range("L10").Value = 1.0
$ws.Range("M10").Value = 2.066403333333334
$ws.Range("N10").Value = 6.199210000000001
$ws.Range("O10").Value = 0.8067501316992219
$ws.Range("P10").Value = 0.8067501316992219
$ws.Range("Q10").Value = 3.721976754353334
$ws.Range("R10").Value = 33.49779078918001
$ws.Range("S10").Value = 0.0897639180685345
$ws.Range("T10").Value = 0.08976391806853448

